# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F) and "最低票价" (G) figures for several
# events on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 7695
    $ws.Range("G2").Value = 138

    $ws.Range("F3").Value = 298

    $ws.Range("F6").Value = 4425

    $ws.Range("F7").Value = 332

    $ws.Range("F9").Value = 281
}
